# BugTriaging-ReadMe.docx edits
#  1. "Dataset Citation :" -> "Dataset Citation:" (drop the gramStart/gramEnd
#     proofErr wrapper Word had inserted around "Citation :" and remove the
#     stray space before the colon).
#  2. The (until-now empty) "Code inspiration:" paragraph loses the explicit
#     paragraph-mark run formatting that had been copied into its <w:pPr>,
#     and gains a trailing plain-text space run.
#  3. The "Code inspiration" hyperlink URL/display text is swapped for the
#     new Kaggle link.
#  4. "in to" -> "into" (again dropping the gramStart/gramEnd proofErr
#     wrapper Word had placed around the misspelling).

$d = $word.ActiveDocument

$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------
# 1. "Dataset Citation :" -> "Dataset Citation:"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Dataset Citation", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$para = $r.Paragraphs(1)
$pr = $para.Range
$attrs = "w14:paraId=`"$($para.Range.Information(10))`""  # placeholder, unused
$xml1 = "<w:p $wNs $w14Ns w14:paraId=`"087131EB`" w14:textId=`"0F1F9D97`" " + `
        "w:rsidR=`"00C5013F`" w:rsidRDefault=`"00C5013F`" w:rsidP=`"002E7EFC`">" + `
        "<w:r><w:t xml:space=`"preserve`">Dataset </w:t></w:r>" + `
        "<w:r><w:t>Citation:</w:t></w:r>" + `
        "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>"
$pr.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2. "Code inspiration:" paragraph - drop paragraph-mark formatting and
#    add a trailing space run.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Code inspiration:", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$para = $r.Paragraphs(1)
$pr = $para.Range
$xml2 = "<w:p $wNs $w14Ns w14:paraId=`"7F7FEBF2`" w14:textId=`"5DC137D0`" " + `
        "w:rsidR=`"00680CDB`" w:rsidRDefault=`"00680CDB`" w:rsidP=`"002E7EFC`">" + `
        "<w:r><w:rPr><w:rFonts w:ascii=`"Helvetica`" w:hAnsi=`"Helvetica`" w:cs=`"Helvetica`"/>" + `
        "<w:color w:val=`"333333`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/>" + `
        "<w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"ECF6FB`"/></w:rPr>" + `
        "<w:t>Code inspiration:</w:t></w:r>" + `
        "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>"
$pr.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3. Swap the "Code inspiration" hyperlink target/display text.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("https://www.machinelearningplus.com/nlp/topic-modeling-gensim-python/", `
                 $true, $false, $false, $false, $false, $true, 1, $false, `
                 "https://www.kaggle.com/selener/multi-class-text-classification-tfidf", 2) | Out-Null

# ---------------------------------------------------------------------
# 4. "in to" -> "into"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("concatenate multiple issue rows", $true, $false, $false, $false, `
                 $false, $true, 1, $false, "", 0) | Out-Null
$para = $r.Paragraphs(1)
$pr = $para.Range
$xml4 = "<w:p $wNs $w14Ns w14:paraId=`"1D46505C`" w14:textId=`"77777777`" " + `
        "w:rsidR=`"00DE5065`" w:rsidRDefault=`"00DE5065`" w:rsidP=`"00DE5065`">" + `
        "<w:r><w:t xml:space=`"preserve`">3. concatenate multiple issue rows </w:t></w:r>" + `
        "<w:r><w:t>into</w:t></w:r>" + `
        "<w:r><w:t xml:space=`"preserve`"> one based on issue id</w:t></w:r></w:p>"
$pr.InsertXML($xml4)
